# Updated Code to Produce Graphical Representation
# - QueueName "Value" (B2) simplified from "API;CustomerSecurityHash" to "API"
# - ReportBackDay "Value" (B8) changed from 1 to 100, and left-aligned

$xlLeft = -4131

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# QueueName row: Value column now just "API"
$ws.Range("B2").Value = "API"

# ReportBackDay row: Value column now 100, left-aligned
$ws.Range("B8").Value = 100
$ws.Range("B8").HorizontalAlignment = $xlLeft
